$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'60.517.35"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -2.47%  "
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'2.899.63"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -3.65%  "
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.11%  "
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'586.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.12%  "
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'147.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.27%  "
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("E7").Value = "'  -0.03%  "
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("D8").Value = "'0.505"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -2.60%  "
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("D9").Value = "'2.898.58"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -3.69%  "
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = "'6.69"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +5.81%  "
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("D11").Value = "'0.145"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -2.93%  "
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("D12").Value = "'0.447"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -2.42%  "
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("E13").Value = "'  -3.04%  "
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "'34.19"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.49%  "
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("E15").Value = "'  +0.39%  "
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = "'3.379.76"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -3.70%  "
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("D17").Value = "'6.81"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -2.26%  "
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'60.457.83"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -2.65%  "
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = "'2.898.49"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -3.74%  "
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'425.59"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -4.71%  "
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21").Value = "'13.63"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -3.83%  "
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = "'0.671"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -2.42%  "
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'7.12"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -3.92%  "
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'80.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -1.79%  "
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "'11.04"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +1.42%  "
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("D26").Value = "'2.19"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -2.04%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("D27").Value = "'11.80"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -2.25%  "
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("E28").Value = "'  +0.03%  "
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("B29").Value = "'NEARProtocol"
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").Value = "'7.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +1.69%  "
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("B30").Value = "'FirstDigitalUSD"
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").Value = "'1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.19%  "
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("B31").Value = "'ImmutableX"
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("C31").Style = "Normal"
$ws.Range("D31").Value = "'2.18"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +3.60%  "
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("B32").Value = "'PancakeSwap"
$ws.Range("B32").Style = "Normal"
$ws.Range("C32").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").Value = "'2.62"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -2.99%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("D33").Value = "'26.51"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -3.36%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("E34").Value = "'  -3.47%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("D35").Value = "'0.0₃0837"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -1.30%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("D36").Value = "'1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -1.88%  "
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("D37").Value = "'5.67"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -2.77%  "
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("B38").Value = "'dogwifhat"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'2.97"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.71%  "
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("D39").Value = "'2.03"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.28%  "
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("B40").Value = "'OKB"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'49.28"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.82%  "
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("D41").Value = "'8.74"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -3.66%  "
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("E42").Value = "'  -1.53%  "
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("D43").Value = "'0.290"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +1.79%  "
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("D44").Value = "'41.70"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.36%  "
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("D45").Value = "'0.0345"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.78%  "
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("D46").Value = "'371.55"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -5.91%  "
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("D47").Value = "'133.23"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.73%  "
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("D48").Value = "'2.652.06"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -2.70%  "
$ws.Range("E48").Style = "Normal"

# Row 50
$ws.Range("D50").Value = "'25.06"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +5.62%  "
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("E51").Value = "'  -1.09%  "
$ws.Range("E51").Style = "Normal"
